$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.303.88"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.08%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.516.55"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.60%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "597.00"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.25%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "173.89"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.65%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.594"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.05%  "

$ws.Range("E9").Value = "  +6.50%  "

$ws.Range("E10").Value = "  -0.35%  "

$ws.Range("E11").Value = "  +0.48%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.126.66"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.67%  "

$ws.Range("E13").Value = "  +0.15%  "

$ws.Range("E14").Value = "  +4.14%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "67.226.77"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.91%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000180"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.75%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.591.10"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.47%  "

$ws.Range("E18").Value = "  +0.31%  "

$ws.Range("E19").Value = "  +2.29%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "396.04"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.07%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "8.03"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.93%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.21"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.26%  "

$ws.Range("E23").Value = "  +0.11%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.539"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.59%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000122"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.71%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.22"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.01%  "

$ws.Range("E27").Value = "  +0.62%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.997"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.28%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.32"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.50%  "

$ws.Range("E30").Value = "  -0.63%  "

$ws.Range("E31").Value = "  +0.95%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "23.94"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.81%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.39"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.59%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.67"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.79%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "163.53"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.60%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.890"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.14%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.92"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.44%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "7.07"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +6.15%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0754"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.26%  "

$ws.Range("E40").Value = "  +0.72%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "26.68"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.60%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "27.36"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.75%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.841.62"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.84%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.62"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.95%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "43.04"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.12%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0305"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.48%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "341.68"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.30%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.09"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.74%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "34.40"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.43%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.51"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.78%  "

$ws.Range("E51").Value = "  -0.70%  "
